$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.622.55'
$ws.Range("E2").Value = '  +0.71%  '

# Row 3
$ws.Range("D3").Value = '3.596.83'
$ws.Range("E3").Value = '  +1.01%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.35%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.30%  '

# Row 9
$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.66%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.136'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '

# Row 11
$ws.Range("E11").Value = '  +0.63%  '

# Row 12
$ws.Range("D12").Value = '4.205.39'
$ws.Range("E12").Value = '  +0.93%  '

# Row 13
$ws.Range("E13").Value = '  +0.73%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.90%  '

# Row 15
$ws.Range("D15").Value = '3.580.93'
$ws.Range("E15").Value = '  +0.57%  '

# Row 16
$ws.Range("D16").Value = '66.697.45'

# Row 17
$ws.Range("E17").Value = '  +0.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '

# Row 19
$ws.Range("E19").Value = '  +1.99%  '

# Row 20
$ws.Range("E20").Value = '  +1.39%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.618'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '78.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").Value = '3.737.62'
$ws.Range("E24").Value = '  +0.93%  '

# Row 25
$ws.Range("E25").Value = '  +0.05%  '

# Row 26
$ws.Range("E26").Value = '  +3.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.20%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.00%  '

# Row 29
$ws.Range("E29").Value = '  -0.28%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("B31").Value = 'RenzoRestakedETH'
$ws.Range("C31").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D31").Value = '3.592.89'
$ws.Range("E31").Value = '  +1.05%  '

# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.159'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.26%  '

# Row 34
$ws.Range("E34").Value = '  -0.25%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.81%  '

# Row 37
$ws.Range("E37").Value = '  +0.36%  '

# Row 38
$ws.Range("E38").Value = '  -2.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '177.75'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.56%  '

# Row 40
$ws.Range("E40").Value = '  +0.23%  '

# Row 41
$ws.Range("E41").Value = '  +0.30%  '

# Row 42
$ws.Range("E42").Value = '  +0.21%  '

# Row 43
$ws.Range("E43").Value = '  -1.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.90%  '

# Row 45
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("E46").Value = '  -1.73%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.64%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.25%  '

# Row 49
$ws.Range("E49").Value = '  +0.56%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.951'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.14%  '

# Row 51
$ws.Range("E51").Value = '  -1.36%  '
